# Generate Report for Handoff
# Updates the localization-status workbook: the two files that were
# previously "Handed back: in sync with en-US" (d4065bc2-... and
# f8802268-...) are now ready for a fresh handoff, with refreshed
# handoff timestamps and a note that the last handback was stale.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet (columns: A File Name, B Path And Name, C Extension,
# D Publish URL, E zh-cn, F de-de, G Latest HO Xliff Generate Date)
# Row 4 -> d4065bc2-f279-4ea1-b0a7-51bd76e4b1a7.md
# Row 5 -> f8802268-b7cf-45a8-b04a-ef668306a399.md
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = $statusReady
$wsOverview.Range("F4").Value = $statusReady
$wsOverview.Range("G4").Value = "2016-09-03 18:37:43"

$wsOverview.Range("E5").Value = $statusReady
$wsOverview.Range("F5").Value = $statusReady
$wsOverview.Range("G5").Value = "2016-09-03 18:37:43"

# ---------------------------------------------------------------------
# zh-cn sheet (columns: A Source File Name, B File Extension, C Status,
# D Source Path, E Priority, F Content Duplicate, G Latest Handoff File,
# H Latest Handoff Datetime, I Latest Target File, J Latest Handback
# File, K Latest Handback DateTime, L Reference Tokens, M To be
# localized, N Dependency From, O Has metadata, P Error Detail)
# Row 4 -> d4065bc2-f279-4ea1-b0a7-51bd76e4b1a7.md
# Row 5 -> f8802268-b7cf-45a8-b04a-ef668306a399.md
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $statusReady
$wsZhCn.Range("H4").Value = "2016-09-03 18:37:38"
$wsZhCn.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ce572c768effbb2b5e36f3fd049013a4f7cdc22/e2e/d4065bc2-f279-4ea1-b0a7-51bd76e4b1a7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb78cd2676906ae67f9cb6e2d6b55b161a56a5a1/e2e/d4065bc2-f279-4ea1-b0a7-51bd76e4b1a7.md."

$wsZhCn.Range("C5").Value = $statusReady
$wsZhCn.Range("H5").Value = "2016-09-03 18:37:38"
$wsZhCn.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ce572c768effbb2b5e36f3fd049013a4f7cdc22/e2e/f8802268-b7cf-45a8-b04a-ef668306a399.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb78cd2676906ae67f9cb6e2d6b55b161a56a5a1/e2e/f8802268-b7cf-45a8-b04a-ef668306a399.md."

# The Error Detail column (P, the 16th column) is now wide enough to
# show the long message without truncation.
$wsZhCn.Columns("P").ColumnWidth = 39.1640625

# ---------------------------------------------------------------------
# de-de sheet (same column layout as zh-cn)
# Row 4 -> d4065bc2-f279-4ea1-b0a7-51bd76e4b1a7.md
# Row 5 -> f8802268-b7cf-45a8-b04a-ef668306a399.md
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $statusReady
$wsDeDe.Range("H4").Value = "2016-09-03 18:37:43"
$wsDeDe.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ce572c768effbb2b5e36f3fd049013a4f7cdc22/e2e/d4065bc2-f279-4ea1-b0a7-51bd76e4b1a7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb78cd2676906ae67f9cb6e2d6b55b161a56a5a1/e2e/d4065bc2-f279-4ea1-b0a7-51bd76e4b1a7.md."

$wsDeDe.Range("C5").Value = $statusReady
$wsDeDe.Range("H5").Value = "2016-09-03 18:37:43"
$wsDeDe.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ce572c768effbb2b5e36f3fd049013a4f7cdc22/e2e/f8802268-b7cf-45a8-b04a-ef668306a399.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb78cd2676906ae67f9cb6e2d6b55b161a56a5a1/e2e/f8802268-b7cf-45a8-b04a-ef668306a399.md."

$wsDeDe.Columns("P").ColumnWidth = 39.1640625
